# Mehr gemacht mit stosse
# For each of the 4 worksheets, insert 7 new columns before column S
# (which pushes the existing "vsp..dvsp_err" block from S:X to Z:AE),
# and populate the freshly inserted S:X columns with the new
# "v_rel / vs_rel / eta" calculations.

$wb = $excel.ActiveWorkbook

$sheetNames = @("Elastisch1", "Elastisch2", "Inelastisch1", "Inelastisch2")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Insert 7 new columns in front of the old "vsp" column (S).
    # The 6 old columns S:X shift to Z:AE, and column Y is left blank.
    $ws.Range("S1:Y1").EntireColumn.Insert()

    # New headers - set in this specific order (S, U, V, T, W, X) to match
    # the order the unique strings were added to the shared string table.
    $ws.Range("S1").Value2 = "v_rel"
    $ws.Range("U1").Value2 = "vs_rel"
    $ws.Range("V1").Value2 = "vs_rel_err"
    $ws.Range("T1").Value2 = "v_rel_err"
    $ws.Range("W1").Value2 = "eta"
    $ws.Range("X1").Value2 = "eta_err"

    # Row 2 formulas (first data row, not part of the shared-formula fill).
    $ws.Range("S2").Formula = "=G2-M2"
    $ws.Range("T2").Formula = "=SQRT(H2^2+N2^2)"
    $ws.Range("U2").Formula = "=J2-P2"
    $ws.Range("V2").Formula = "=SQRT(K2^2+Q2^2)"
    $ws.Range("W2").Formula = "=U2^2/S2^2"
    $ws.Range("X2").Formula = "=2*ABS(U2/S2)*SQRT((V2/S2)^2+(U2*T2/S2^2)^2)"

    # Rows 3:6 - filled down together so they form shared formula groups.
    $ws.Range("S3:S6").Formula = "=G3-M3"
    $ws.Range("T3:T6").Formula = "=SQRT(H3^2+N3^2)"
    $ws.Range("U3:U6").Formula = "=J3-P3"
    $ws.Range("V3:V6").Formula = "=SQRT(K3^2+Q3^2)"
    $ws.Range("W3:W6").Formula = "=U3^2/S3^2"
    $ws.Range("X3:X6").Formula = "=2*ABS(U3/S3)*SQRT((V3/S3)^2+(U3*T3/S3^2)^2)"
}

# Restore the active-cell selections that were recorded for each sheet.
$wb.Worksheets.Item("Elastisch1").Range("W4").Select()
$wb.Worksheets.Item("Elastisch2").Range("N6").Select()
$wb.Worksheets.Item("Inelastisch1").Range("F18").Select()
$wb.Worksheets.Item("Inelastisch2").Range("R1").Select()

# Inelastisch2 (sheet4 / tab index 4) remains the active tab.
$wb.Worksheets.Item("Inelastisch2").Activate()
